# control de coma por punto decimal
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 through 10 (the table is shrinking from A1:D10 to A1:D5)
$ws.Range("A6:D10").EntireRow.Delete()

# Update the remaining data rows (2-5) with the new values
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = -250
$ws.Range("D2").Value = 1.000005

$ws.Range("B3").Value = 40
$ws.Range("C3").Value = -1.20892581961463 * [Math]::Pow(10, 24)
$ws.Range("D3").Value = 1.000005

$ws.Range("B4").Value = 4
$ws.Range("C4").Value = -250
$ws.Range("D4").Value = 36

$ws.Range("B5").Value = 4
$ws.Range("C5").Value = -250
$ws.Range("D5").Value = 0
